$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new nombre_aides C, new montant_total E)
$updates = @{
    3   = @(249327, 1036476130)
    4   = @(103466, 627499207)
    62  = @(4189, 9186399)
    91  = @(151105, 482131030)
    92  = @(409019, 1593716220)
    93  = @(209507, 1307872027)
    94  = @(94148, 915830612)
    95  = @(50728, 930867754)
    96  = @(17247, 789598588)
    104 = @(135223, 272132060)
    114 = @(3799, 9102060)
    115 = @(11692, 32953239)
    116 = @(4560, 20567318)
    127 = @(57, 207876)
    131 = @(75581, 307211165)
    157 = @(12097, 182911050)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
